$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Delete review row 10 (hermanliran@gmail.com / nevilgreen@gmail.com entry)
#    This shifts rows 11-17 up to become rows 10-16.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Delete()

# ---------------------------------------------------------------------------
# 2. Append a new review as row 17, copying the cell formatting from row 14
#    (a "taxi" row with the same style pattern) and then filling in values.
# ---------------------------------------------------------------------------
$cols = @("A","B","C","D","E","F","G")
foreach ($col in $cols) {
    $ws.Range($col + "14").Copy()
    $ws.Range($col + "17").PasteSpecial(-4122)
}

$ws.Range("A17").Value = "com.singleton.strechy"
$ws.Range("B17").Value = "taxi"
$ws.Range("C17").Value = "georggini2@gmail.com"
$ws.Range("D17").Value = "jorjkluni03@gmail.com"
$ws.Range("E17").Value = "27/5/2019 15:59"
$ws.Range("F17").Value = "one of the greatest assets in my phone device. Guaranteed!"
$ws.Range("G17").Value = "confirm"

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks collection to reflect the new row numbering and
#    the newly added reviewer e-mail hyperlink.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:nevilgreen12@gmail.com", "", "", "nevilgreen12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:rocketaso@gmail.com", "", "", "rocketaso@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:danfogel100@gmail.com", "", "", "danfogel100@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:avishaybar12@gmail.com", "", "", "avishaybar12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:danfogel100@gmail.com", "", "", "danfogel100@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:avishaybar12@gmail.com", "", "", "avishaybar12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com")

# ---------------------------------------------------------------------------
# 4. Restore the selected cell to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B10").Select()
